$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.223.04"
$ws.Cells.Item(2, 5).Value = "  +3.10%  "
$ws.Cells.Item(3, 4).Value = "1.897.46"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).Value = "'325.20"
$ws.Cells.Item(5, 5).Value = "  +3.44%  "
$ws.Cells.Item(6, 5).Value = "  -0.22%  "
$ws.Cells.Item(7, 4).Value = "'0.5159"
$ws.Cells.Item(7, 5).Value = "  +0.22%  "
$ws.Cells.Item(8, 4).Value = "'0.4008"
$ws.Cells.Item(8, 5).Value = "  +2.20%  "
$ws.Cells.Item(9, 4).Value = "'0.08447"
$ws.Cells.Item(9, 5).Value = "  +0.28%  "
$ws.Cells.Item(10, 4).Value = "'42.68"
$ws.Cells.Item(10, 5).Value = "  +0.91%  "
$ws.Cells.Item(11, 5).Value = "  +0.29%  "
$ws.Cells.Item(12, 4).Value = "'23.42"
$ws.Cells.Item(12, 5).Value = "  +13.20%  "
$ws.Cells.Item(13, 4).Value = "'6.428"
$ws.Cells.Item(13, 5).Value = "  +2.54%  "
$ws.Cells.Item(14, 4).Value = "1.891.67"
$ws.Cells.Item(14, 5).Value = "  +0.04%  "
$ws.Cells.Item(15, 4).Value = "'7.338"
$ws.Cells.Item(15, 5).Value = "  +0.86%  "
$ws.Cells.Item(17, 4).Value = "'94.77"
$ws.Cells.Item(17, 5).Value = "  +1.78%  "
$ws.Cells.Item(18, 5).Value = "  +0.56%  "
$ws.Cells.Item(19, 4).Value = "'0.06644"
$ws.Cells.Item(19, 5).Value = "  -1.27%  "
$ws.Cells.Item(20, 5).Value = "  +2.27%  "
$ws.Cells.Item(21, 5).Value = "  -0.29%  "
$ws.Cells.Item(22, 4).Value = "'5.954"
$ws.Cells.Item(22, 5).Value = "  -0.92%  "
$ws.Cells.Item(23, 4).Value = "30.227.69"
$ws.Cells.Item(23, 5).Value = "  +3.06%  "
$ws.Cells.Item(24, 5).Value = "  +1.40%  "
$ws.Cells.Item(25, 5).Value = "  +0.60%  "
$ws.Cells.Item(26, 4).Value = "2.110.60"
$ws.Cells.Item(26, 5).Value = "  +0.12%  "
$ws.Cells.Item(27, 5).Value = "  +3.58%  "
$ws.Cells.Item(28, 4).Value = "'161.08"
$ws.Cells.Item(28, 5).Value = "  +1.31%  "
$ws.Cells.Item(29, 4).Value = "'2.363"
$ws.Cells.Item(29, 5).Value = "  -2.79%  "
$ws.Cells.Item(30, 4).Value = "'128.87"
$ws.Cells.Item(30, 5).Value = "  +1.38%  "
$ws.Cells.Item(31, 5).Value = "  +3.49%  "
$ws.Cells.Item(32, 4).Value = "'0.1057"
$ws.Cells.Item(33, 4).Value = "'6.080"
$ws.Cells.Item(33, 5).Value = "  -1.04%  "
$ws.Cells.Item(34, 4).Value = "'3.761"
$ws.Cells.Item(34, 5).Value = "  +2.74%  "
$ws.Cells.Item(35, 4).Value = "'0.02494"
$ws.Cells.Item(35, 5).Value = "  +0.51%  "
$ws.Cells.Item(36, 4).Value = "'0.06558"
$ws.Cells.Item(36, 5).Value = "  -0.28%  "
$ws.Cells.Item(37, 4).Value = "'5.275"
$ws.Cells.Item(38, 5).Value = "  +0.24%  "
$ws.Cells.Item(39, 5).Value = "  -0.83%  "
$ws.Cells.Item(40, 4).Value = "'11.78"
$ws.Cells.Item(40, 5).Value = "  +4.62%  "
$ws.Cells.Item(41, 4).Value = "'0.6500"
$ws.Cells.Item(41, 5).Value = "  -0.25%  "
$ws.Cells.Item(42, 4).Value = "'8.704"
$ws.Cells.Item(42, 5).Value = "  -3.42%  "
$ws.Cells.Item(43, 4).Value = "'1.231"
$ws.Cells.Item(43, 5).Value = "  -0.19%  "
$ws.Cells.Item(44, 5).Value = "  +0.75%  "
$ws.Cells.Item(45, 4).Value = "'13.26"
$ws.Cells.Item(45, 5).Value = "  +0.15%  "
$ws.Cells.Item(46, 4).Value = "'3.702"
$ws.Cells.Item(46, 5).Value = "  +0.81%  "
$ws.Cells.Item(47, 4).Value = "'2.058"
$ws.Cells.Item(47, 5).Value = "  +0.33%  "
$ws.Cells.Item(48, 4).Value = "'1.236"
$ws.Cells.Item(48, 5).Value = "  +0.50%  "
$ws.Cells.Item(49, 4).Value = "'124.42"
$ws.Cells.Item(49, 5).Value = "  +0.71%  "
$ws.Cells.Item(50, 4).Value = "'1.165"
$ws.Cells.Item(50, 5).Value = "  +0.74%  "
$ws.Cells.Item(51, 5).Value = "  +1.74%  "
